$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 10 weekly "Poroto granado" records in rows 43-52.
# Two new weekly records need to be inserted into this block (one becoming the
# new row 43, another becoming the new row 46), pushing the existing records
# down so the table ends at row 54 instead of row 52.

# Insert a new blank row at row 43 - this pushes the current rows 43-52 down to 44-53.
$ws.Rows.Item(43).Insert()

# Insert a second new blank row at row 46 - this pushes rows 46-53 down to 47-54.
$ws.Rows.Item(46).Insert()

# Fill in the first newly inserted row (row 43) with its data.
$ws.Cells.Item(43, 1).Value = 4
$ws.Cells.Item(43, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(43, 3).Value = "Los Lagos"
$ws.Cells.Item(43, 4).Value = 44908
$ws.Cells.Item(43, 5).Value = 10
$ws.Cells.Item(43, 6).Value = 100112030
$ws.Cells.Item(43, 7).Value = "Poroto granado"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 60
$ws.Cells.Item(43, 11).Value = 50000
$ws.Cells.Item(43, 12).Value = 50000
$ws.Cells.Item(43, 13).Value = 50000
$ws.Cells.Item(43, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(43, 15).Value = "Región Metropolitana"
$ws.Cells.Item(43, 16).Value = 2000
$ws.Cells.Item(43, 17).Value = 25
$ws.Cells.Item(43, 18).Value = "Hortaliza"

# Fill in the second newly inserted row (row 46) with its data.
$ws.Cells.Item(46, 1).Value = 4
$ws.Cells.Item(46, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(46, 3).Value = "Los Lagos"
$ws.Cells.Item(46, 4).Value = 44932
$ws.Cells.Item(46, 5).Value = 10
$ws.Cells.Item(46, 6).Value = 100112030
$ws.Cells.Item(46, 7).Value = "Poroto granado"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 80
$ws.Cells.Item(46, 11).Value = 47000
$ws.Cells.Item(46, 12).Value = 47000
$ws.Cells.Item(46, 13).Value = 47000
$ws.Cells.Item(46, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(46, 15).Value = "Región Metropolitana"
$ws.Cells.Item(46, 16).Value = 1880
$ws.Cells.Item(46, 17).Value = 25
$ws.Cells.Item(46, 18).Value = "Hortaliza"
